$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.800.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.309.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.81%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.79'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.33'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.22'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.60'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.11%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.976'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.47%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.660.03'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.310.46'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.766.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.25%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.39'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +32.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.86'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.55'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.39'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.28%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.68'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.16'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +7.29%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.28'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +7.67%  '

$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.90'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0893'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.132'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.59'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -8.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.115'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.62'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0354'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.71'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.73'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.58'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +8.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.63'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.53'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.226'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.96%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.43'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '82.16'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '114.68'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.31'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.90'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.619.03'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.74%  '
